# Issue #73: Standardize title-cased values to lowercase (#84)
#
# The "Sex" and "Age" columns contain a mix of title-cased category labels
# (Male/Female, Adult/Subadult/Fawn). Standardize them to lowercase across
# every sheet in the workbook.

$wb = $excel.ActiveWorkbook

$map = @{
    "Male"     = "male"
    "Female"   = "female"
    "Adult"    = "adult"
    "Subadult" = "subadult"
    "Fawn"     = "fawn"
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $current = $cell.Value2
            if ($null -ne $current -and $map.ContainsKey($current)) {
                $cell.Value = $map[$current]
            }
        }
    }
}
